# Generate Report for Handoff
# Updates the localization-status workbook after a handoff run: renames the
# source file's generated guid, refreshes timestamps, updates the handoff
# xliff hashes, and clears the (not-yet-handed-back) target/handback columns.

$wb = $excel.ActiveWorkbook

$oldGuid = "6d63531b-c377-4854-bb56-6b0546cfa32e"
$newGuid = "782ec4a5-2018-4d63-bfc3-9a4fc9fa655d"

$oldHash = "fd6b7621bac04a3b6b2391af348f717d54342167"
$newHash = "a9ac18e6b7672971e1b61c69c2e1edff0dad2cee"

# ---- Overview sheet ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("A2").Value = "$newGuid.md"
$ov.Range("B2").Value = "e2e\$newGuid.md"
$ov.Hyperlinks.Item(1).TextToDisplay = "e2e\$newGuid.md"
$ov.Range("G2").Value = "2016-09-05 13:14:16"

# ---- zh-cn sheet ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("A2").Value = "$newGuid.md"
$zh.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
$zh.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$zh.Range("H2").Value = "2016-09-05 13:14:10"
$zh.Range("I2").Value = ""
$zh.Range("J2").Value = ""
$zh.Range("K2").Value = "0001-01-01 00:00:00"
$zh.Hyperlinks.Item(2).Delete()

# ---- de-de sheet ----
$de = $wb.Worksheets.Item("de-de")
$de.Range("A2").Value = "$newGuid.md"
$de.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
$de.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$de.Range("I2").Value = ""
$de.Range("J2").Value = ""
$de.Range("K2").Value = "0001-01-01 00:00:00"
$de.Hyperlinks.Item(2).Delete()

# ---- column widths (narrowed now that Latest Target File / Latest Handback
#      File are empty) ----
$zh.Columns.Item(9).ColumnWidth = 18.6506053379604
$zh.Columns.Item(10).ColumnWidth = 21.7054770333426
$de.Columns.Item(9).ColumnWidth = 18.6506053379604
$de.Columns.Item(10).ColumnWidth = 21.7054770333426
